# Refresh the cryptos list (prices + 1h volume%) per the GitHub Actions commit.
# All D/E columns are plain text cells (t="inlineStr") holding formatted strings
# (e.g. "529.84", "  +2.19%  "). Some new price strings parse as valid numbers,
# so for those we briefly force Text format before assigning, then restore the
# default "Normal" style so no stray number-format styling is left behind -
# otherwise Excel's autodetect would convert the cell to a Number and could
# even drop a significant trailing zero (e.g. "22.30" -> 22.3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.608.90"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "3.163.46"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.87%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.442"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.35%  "
$ws.Range("E11").Value = "  +4.90%  "
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("D13").Value = "3.712.17"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.71%  "
$ws.Range("D16").Value = "58.653.75"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.07%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.164.17"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.533"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.06%  "
$ws.Range("D29").Value = "0.0₃0873"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("E35").Value = "  +4.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  +5.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0699"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").Value = "2.661.07"
$ws.Range("E41").Value = "  +7.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.723"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("E45").Value = "  +8.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "3.207.25"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.49%  "
